$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-30 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-31 Saturday", 2)
$d.Content.Find.Execute("32×88=2816", $true, $false, $false, $false, $false, $true, 1, $false, "56×96=5376", 2)
$d.Content.Find.Execute("52×35=1820", $true, $false, $false, $false, $false, $true, 1, $false, "31×86=2666", 2)
$d.Content.Find.Execute("97×76=7372", $true, $false, $false, $false, $false, $true, 1, $false, "13×28=364", 2)
$d.Content.Find.Execute("31×55=1705", $true, $false, $false, $false, $false, $true, 1, $false, "94×71=6674", 2)
$d.Content.Find.Execute("52×78=4056", $true, $false, $false, $false, $false, $true, 1, $false, "52×50=2600", 2)
$d.Content.Find.Execute("79×95=7505", $true, $false, $false, $false, $false, $true, 1, $false, "34×66=2244", 2)
$d.Content.Find.Execute("33×51=1683", $true, $false, $false, $false, $false, $true, 1, $false, "67×40=2680", 2)
$d.Content.Find.Execute("53×52=2756", $true, $false, $false, $false, $false, $true, 1, $false, "46×90=4140", 2)
$d.Content.Find.Execute("91×47=4277", $true, $false, $false, $false, $false, $true, 1, $false, "95×54=5130", 2)
$d.Content.Find.Execute("46×99=4554", $true, $false, $false, $false, $false, $true, 1, $false, "67×37=2479", 2)
$d.Content.Find.Execute("75×86=6450", $true, $false, $false, $false, $false, $true, 1, $false, "89×98=8722", 2)
$d.Content.Find.Execute("95×33=3135", $true, $false, $false, $false, $false, $true, 1, $false, "35×35=1225", 2)
$d.Content.Find.Execute("68×82=5576", $true, $false, $false, $false, $false, $true, 1, $false, "99×34=3366", 2)
$d.Content.Find.Execute("92×66=6072", $true, $false, $false, $false, $false, $true, 1, $false, "80×80=6400", 2)
$d.Content.Find.Execute("26×37=962", $true, $false, $false, $false, $false, $true, 1, $false, "69×84=5796", 2)
$d.Content.Find.Execute("79×55=4345", $true, $false, $false, $false, $false, $true, 1, $false, "55×49=2695", 2)
$d.Content.Find.Execute("92×32=2944", $true, $false, $false, $false, $false, $true, 1, $false, "99×59=5841", 2)
$d.Content.Find.Execute("22×95=2090", $true, $false, $false, $false, $false, $true, 1, $false, "99×51=5049", 2)
$d.Content.Find.Execute("40×83=3320", $true, $false, $false, $false, $false, $true, 1, $false, "45×24=1080", 2)
$d.Content.Find.Execute("67×45=3015", $true, $false, $false, $false, $false, $true, 1, $false, "87×40=3480", 2)
$d.Content.Find.Execute("35×43=1505", $true, $false, $false, $false, $false, $true, 1, $false, "82×41=3362", 2)
$d.Content.Find.Execute("42×17=714", $true, $false, $false, $false, $false, $true, 1, $false, "40×50=2000", 2)
$d.Content.Find.Execute("26×48=1248", $true, $false, $false, $false, $false, $true, 1, $false, "64×14=896", 2)
$d.Content.Find.Execute("47×71=3337", $true, $false, $false, $false, $false, $true, 1, $false, "50×11=550", 2)
$d.Content.Find.Execute("83×73=6059", $true, $false, $false, $false, $false, $true, 1, $false, "41×11=451", 2)
